$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H ("15-jun") mirrors the formatting of the existing G column
# (header style + numeric/centered body style), so copy formats across first.
$ws.Range("G1:G11").Copy()
$ws.Range("H1:H11").PasteSpecial(-4122)

# Header for the new date column
$ws.Range("H1").Value = "15-jun"

# Daily counts for the new date
$ws.Range("H2").Value = 15
$ws.Range("H3").Value = 13
$ws.Range("H4").Value = 7
$ws.Range("H5").Value = 13
$ws.Range("H6").Value = 17
$ws.Range("H7").Value = 13
$ws.Range("H8").Value = 10
$ws.Range("H9").Value = 15
$ws.Range("H10").Value = 18
$ws.Range("H11").Value = 8

# Selection moved from G12 to F12 in the saved view
$ws.Range("F12").Select()
